$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Columns("A:A").Insert()

# New column A: orientation header + value
$ws.Range("A1").Value = "orientation"
$ws.Range("A2").Value = "[-90, -45, 0, 45, 90]"
$ws.Columns("A:A").ColumnWidth = 16.83

# Update trial_duration value (now in column D) from 5 to 2
$ws.Range("D2").Value = 2

# Update isi value (now in column F) - set to 1
$ws.Range("F2").Value = 1

# Update selection to match target state
$null = $ws.Range("F2").Select()
